$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.248.86"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.511.23"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "538.32"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").Value = "137.82"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "2.524.70"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("D14").Value = "2.960.52"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "23.17"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "59.138.05"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "2.528.49"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "11.13"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "325.91"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "65.52"
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  -1.30%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "6.73"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  +6.36%  "
$ws.Range("D33").Value = "165.10"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.47"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("D36").Value = "18.49"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "4.13"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "36.80"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "290.04"
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "131.76"
$ws.Range("E45").Value = "  +8.42%  "
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "0.0933"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "17.42"
$ws.Range("E51").Value = "  -1.26%  "
